# Update the "Förändrad" (Changed) date column (C) for rows 2-13 from
# 2023-10-05 (serial 45204) to 2023-10-08 (serial 45207).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C13").Value = 45207
